$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (M2:T2)
$ws.Range("M2").Value = 0.8155003333333334
$ws.Range("N2").Value = 2.446501
$ws.Range("O2").Value = 0.1910612426590028
$ws.Range("P2").Value = 0.1910612426590029
$ws.Range("Q2").Value = 0.05886036755900001
$ws.Range("R2").Value = 0.529743308031
$ws.Range("S2").Value = 0.1910612426590028
$ws.Range("T2").Value = 0.1910612426590029

# Update row 3 values (O3,P3,S3,T3)
$ws.Range("O3").Value = 0.7809105179307759
$ws.Range("P3").Value = 0.780910517930776
$ws.Range("S3").Value = 0.7809105179307759
$ws.Range("T3").Value = 0.780910517930776

# Update row 4 values (M4:T4)
$ws.Range("M4").Value = 0.119632
$ws.Range("N4").Value = 0.358896
$ws.Range("O4").Value = 0.02802823941022116
$ws.Range("P4").Value = 0.02802823941022117
$ws.Range("Q4").Value = 0.008634678864000002
$ws.Range("R4").Value = 0.077712109776
$ws.Range("S4").Value = 0.02802823941022116
$ws.Range("T4").Value = 0.02802823941022117

# Delete row 5 entirely
$ws.Rows("5:5").Delete()
